# Add a new "footertext" worksheet (config data used by the regression
# suite / test-case description, per the commit message) with
# address / email / phone columns, after the existing "subcategory" sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "footertext"

# Header row
$ws.Range("A1").Value = "address"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "phone"

# Data row
$ws.Range("A2").Value = "Mr John Smith`n132 My Street, `nKingston, `nNew York 12401."
$ws.Range("B2").Value = "Johnsmith"
$ws.Range("C2").Value = "abcdefgh++++"

# Match the look of the workbook's other data sheets (e.g. "subcategory"):
# copy its cell formatting (font/alignment) onto the new sheet's cells.
$subSheet = $wb.Worksheets.Item("subcategory")
$subSheet.Range("A1").Copy()
$ws.Range("A1:C2").PasteSpecial(-4122)  # xlPasteFormats

# The multi-line address needs to wrap within its cell.
$ws.Range("A2").WrapText = $true

# Restore the original active sheet/selection.
$wb.Worksheets.Item("loginpage").Activate()
